# Replace the "JSONATA" column (G) + separate "Mapping Result" column (H)
# on the "TS Parameters" sheet with a single "Mapping Result" column (G).
#   - G1 header: "JSONATA" -> "Mapping Result"
#   - Where a result already existed in H, move that value into G
#     (replacing the JSONATA formula text) and remove H.
#   - Where G held a JSONATA formula but H had no result, blank G to a
#     single space " ".
#   - Where G was previously empty, add a new G cell containing " ".
#   - Rows 21, 26 and 31 are left untouched (JSONATA text stays in G).
#   - Column H is cleared entirely, shrinking the used range to A1:G60.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TS Parameters")

# 1) Header rename
$ws.Range("G1").Value = "Mapping Result"

# 2) Cells whose H-column result replaces the G-column JSONATA formula.
#    G2/G3 originally use the "quote prefix" style (s="34"); a leading
#    apostrophe keeps that formatting (forces text) exactly as plain
#    Excel data-entry would, while the apostrophe itself is not stored.
$bullet = [char]0x2022
$g30Text = $bullet + "`t" + "To evaluate the effect of intravenous infusions of Beroclin administered once daily on motor symptoms in subjects with early stage Parkinson s disease."

$ws.Range("G2").Value = "'Y"
$ws.Range("G3").Value = "'" + '18 ["Year","Year"] '
$ws.Range("G4").Value = '100 ["Year","Year"] '
$ws.Range("G22").Value = "PARALLEL"
$ws.Range("G23").Value = "Pharmacologic Substance"
$ws.Range("G30").Value = $g30Text
$ws.Range("G53").Value = "Double Blind Study"
$ws.Range("G57").Value = "Safety and Efficacy of the Xanomeline Transdermal Therapeutic System (TTS) in Patients with Mild to Moderate Alzheimer's Disease"

# 3) Existing G cells (JSONATA formula, no matching H result) blanked to " "
$blankExistingG = @(5, 6, 7, 29, 41, 49, 58)
foreach ($r in $blankExistingG) {
    $ws.Cells.Item($r, 7).Value = " "
}

# 4) Rows that previously had no G cell at all get a new one containing " "
$newG = @(8,9,10,11,12,13,14,15,16,17,18,19,20,24,25,27,28,32,33,34,35,36,37,38,39,40,42,43,44,45,46,47,48,50,51,52,54,55,56,59,60)
foreach ($r in $newG) {
    $ws.Cells.Item($r, 7).Value = " "
}

# Rows 21, 26 and 31 keep their original JSONATA text in G - no change.

# 5) Remove column H entirely (its values have been folded into G above)
$ws.Range("H1:H60").ClearContents()
